$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("ランサーズ"): prepend 5 newly scraped job postings at the top.
# The existing rows 2-107 shift down to rows 7-112.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Grab the hyperlink-formatted cell first so we can copy its formatting
# (style index) onto the freshly inserted rows without creating new style
# entries in styles.xml.
$ws1.Range("F2").Copy()

# Insert 5 blank rows above the current row 2 (pushes old rows 2-107 -> 7-112).
$ws1.Rows("2:6").Insert()

# Re-apply the hyperlink cell formatting (style only) to the new F2:F6 cells.
$ws1.Range("F2:F6").PasteSpecial(-4122)

# --- Row 2 -----------------------------------------------------------------
$ws1.Range("A2").Value = "2025-09-03 12:34:36"
$ws1.Range("B2").Value = "React / React Native|恋愛×AIアプリのフロントエンジニア募集!急成長プロダクト"
$ws1.Range("C2").Value = "システム開発"
$ws1.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws1.Range("E2").Value = "期限情報なし"
$ws1.Hyperlinks.Add($ws1.Range("F2"), "https://www.lancers.jp/work/detail/5385757")
$ws1.Range("G2").Value = 445
$ws1.Range("H2").Value = "🔥AI,React ◇アプリ"

# --- Row 3 -----------------------------------------------------------------
$ws1.Range("A3").Value = "2025-09-03 12:34:36"
$ws1.Range("B3").Value = "【急募】Teamsチャット履歴をPythonでテキスト出力したい"
$ws1.Range("C3").Value = "システム開発"
$ws1.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws1.Range("E3").Value = "期限情報なし"
$ws1.Hyperlinks.Add($ws1.Range("F3"), "https://www.lancers.jp/work/detail/5385818")
$ws1.Range("G3").Value = 190
$ws1.Range("H3").Value = "🔥Python"

# --- Row 4 -----------------------------------------------------------------
$ws1.Range("A4").Value = "2025-09-03 12:34:36"
$ws1.Range("B4").Value = "【急募】多言語対応SaaSのフロントエンド開発者募集"
$ws1.Range("C4").Value = "システム開発"
$ws1.Range("D4").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws1.Range("E4").Value = "期限情報なし"
$ws1.Hyperlinks.Add($ws1.Range("F4"), "https://www.lancers.jp/work/detail/5385914")
$ws1.Range("G4").Value = 75
$ws1.Range("H4").Value = "◆開発"

# --- Row 5 -----------------------------------------------------------------
$ws1.Range("A5").Value = "2025-09-03 12:34:36"
$ws1.Range("B5").Value = "Contact Form7×メールサーバ×受発注システム連携"
$ws1.Range("C5").Value = "システム開発"
$ws1.Range("D5").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws1.Range("E5").Value = "期限情報なし"
$ws1.Hyperlinks.Add($ws1.Range("F5"), "https://www.lancers.jp/work/detail/5385814")
$ws1.Range("G5").Value = 33

# --- Row 6 -----------------------------------------------------------------
$ws1.Range("A6").Value = "2025-09-03 12:34:36"
$ws1.Range("B6").Value = "限定公開 限定公開の仕事"
$ws1.Range("C6").Value = "システム開発"
$ws1.Range("D6").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws1.Range("E6").Value = "期限情報なし"
$ws1.Hyperlinks.Add($ws1.Range("F6"), "https://www.lancers.jp/work/detail/5385681")
$ws1.Range("G6").Value = 25

# ---------------------------------------------------------------------------
# Sheet 2 ("統計"): append a new stats snapshot row at the bottom (row 59).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A59").Value = "2025-09-03T12:34:36.566706"
$ws2.Range("B59").Value = 21
$ws2.Range("C59").Value = "全案件リスト"
$ws2.Range("D59").Value = 66.7
$ws2.Range("E59").Value = 7
$ws2.Range("F59").Value = 8
$ws2.Range("G59").Value = 21
